$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" - refresh the handoff timestamps for the
# d8eb904b file (zh-cn + de-de + overview rollup) and, in de-de, also the
# b58c571e file which shared the same prior timestamp value.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: row 7 = d8eb904b-9657-4841-a68b-b67d361c457a.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-08-15 18:38:14"

# zh-cn table: row 7 = d8eb904b-9657-4841-a68b-b67d361c457a.md
# Column H = "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-08-15 18:37:59"

# de-de table: row 6 = b58c571e-5a42-4543-8519-35840c454a8a.md
#              row 7 = d8eb904b-9657-4841-a68b-b67d361c457a.md
# Column H = "Latest Handoff Datetime"
$wsDeDe.Range("H6").Value = "2016-08-15 18:38:14"
$wsDeDe.Range("H7").Value = "2016-08-15 18:38:14"
